$d = $word.ActiveDocument
$paras = $d.Paragraphs

# Find the paragraph that holds the "Ver no Jupiter Salvar em pdf Salvar em docx"
# run -- this is the unique anchor for the block that must be removed.
$jupiterIndex = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.StartsWith("Ver no Jupiter Salvar em pdf Salvar em docx")) {
        $jupiterIndex = $i
    }
}

# The block to delete is:
#   jupiterIndex - 1 : empty paragraph right after "LOM3003: ..."
#   jupiterIndex     : "Ver no Jupiter Salvar em pdf Salvar em docx"
#   jupiterIndex + 1 : empty paragraph
#   jupiterIndex + 2 : empty paragraph with pageBreakBefore
# so that the paragraph with "LOM3003: ..." is directly followed by the
# paragraph that used to be jupiterIndex + 3 (an empty Normal paragraph).

$startPara = $paras.Item($jupiterIndex - 1)
$endPara = $paras.Item($jupiterIndex + 2)

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range.Delete()
